# Update countries & provincias Spain
# Applies the daily refresh of the COVID "Pais" dataset:
#  - Country-name relabeling caused by the source re-sorting some rows
#    (Crucero/Sri Lanka swap, a 7-row Nepal rotation, Namibia/San Vicente swap)
#  - Updated case/recovered/death counts for a number of countries

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row label (country name, column A) changes -----------------------
$ws.Range("A103").Value = "Sri Lanka"
$ws.Range("A104").Value = "Crucero"

$ws.Range("A160").Value = "Bahamas"
$ws.Range("A161").Value = "Guyana"
$ws.Range("A162").Value = "Liechtenstein"
$ws.Range("A163").Value = "Barbados"
$ws.Range("A164").Value = "Mozambique"
$ws.Range("A165").Value = "San Martin (Parte Holandesa)"
$ws.Range("A166").Value = "Nepal"

$ws.Range("A194").Value = "Namibia"
$ws.Range("A195").Value = "San Vicente y las Granadinas"

# --- Updated statistics (columns B:H) ----------------------------------
# row, B, C, D, E, F, G, H
$rows = @(
    @(4,   1178574, 17800, 177705, 932655, 16156, 770, 68214),
    @(8,   168693,  297,   50784,  93014,  3819,  135, 24895),
    @(9,   165383,  416,   130600, 27943,  1979,  28,  6840),
    @(20,  29905,   88,    24500,  3643,   141,   0,   1762),
    @(46,  7833,    24,    32,     7590,   37,    0,   211),
    @(62,  3383,    99,    1718,   1657,   1,     0,   8),
    @(103, 718,     16,    184,    527,    1,     0,   7),
    @(104, 712,     0,     645,    54,     4,     0,   13),
    @(112, 563,     19,    213,    323,    0,     1,   27),
    @(120, 461,     1,     367,    85,     5,     0,   9),
    @(160, 83,      0,     24,     48,     1,     0,   11),
    @(161, 82,      0,     22,     51,     2,     0,   9),
    @(162, 82,      0,     55,     26,     0,     0,   1),
    @(163, 81,      0,     44,     30,     4,     0,   7),
    @(164, 79,      0,     18,     61,     0,     0,   0),
    @(165, 76,      0,     44,     19,     7,     0,   13),
    @(166, 75,      16,    16,     59,     0,     0,   0),
    @(208, 10,      0,     9,      0,      0,     0,   1)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("B$rowNum").Value = $r[1]
    $ws.Range("C$rowNum").Value = $r[2]
    $ws.Range("D$rowNum").Value = $r[3]
    $ws.Range("E$rowNum").Value = $r[4]
    $ws.Range("F$rowNum").Value = $r[5]
    $ws.Range("G$rowNum").Value = $r[6]
    $ws.Range("H$rowNum").Value = $r[7]
}
